$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.884.20'
$ws.Range("E2").Value = '  -1.08%  '

# Row 3
$ws.Range("D3").Value = '2.193.33'
$ws.Range("E3").Value = '  -2.37%  '

# Row 4
$ws.Range("E4").Value = '  -0.22%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '294.32'
$c.ClearFormats()
$ws.Range("E5").Value = '  -4.23%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '87.86'
$c.ClearFormats()
$ws.Range("E6").Value = '  -6.34%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.565'
$c.ClearFormats()
$ws.Range("E7").Value = '  -0.94%  '

# Row 8
$ws.Range("E8").Value = '  -0.12%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.477'
$c.ClearFormats()
$ws.Range("E9").Value = '  -9.01%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '31.97'
$c.ClearFormats()
$ws.Range("E10").Value = '  -7.83%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0766'
$c.ClearFormats()
$ws.Range("E11").Value = '  -5.54%  '

# Row 12
$ws.Range("E12").Value = '  -1.66%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '6.72'
$c.ClearFormats()
$ws.Range("E13").Value = '  -6.27%  '

# Row 14
$ws.Range("D14").Value = '2.527.38'
$ws.Range("E14").Value = '  -2.42%  '

# Row 15
$ws.Range("D15").Value = '2.269.42'
$ws.Range("E15").Value = '  -2.91%  '

# Row 16
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.765'
$c.ClearFormats()
$ws.Range("E16").Value = '  -8.88%  '

# Row 17
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '12.91'
$c.ClearFormats()
$ws.Range("E17").Value = '  -5.55%  '

# Row 18
$ws.Range("D18").Value = '43.511.56'
$ws.Range("E18").Value = '  -1.24%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0880'
$ws.Range("E19").Value = '  -8.55%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '5.79'
$c.ClearFormats()
$ws.Range("E20").Value = '  -9.42%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '10.61'
$c.ClearFormats()
$ws.Range("E21").Value = '  -15.77%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '62.65'
$c.ClearFormats()
$ws.Range("E22").Value = '  -4.74%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '229.37'
$c.ClearFormats()
$ws.Range("E23").Value = '  -3.31%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.74'
$c.ClearFormats()
$ws.Range("E24").Value = '  -8.63%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.ClearFormats()
$ws.Range("E25").Value = '  +0.58%  '

# Row 26
$ws.Range("E26").Value = '  -9.03%  '

# Row 27
$ws.Range("E27").Value = '  +0.35%  '

# Row 28
$ws.Range("B28").Value = 'InjectiveProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '35.40'
$c.ClearFormats()
$ws.Range("E28").Value = '  -8.26%  '

# Row 29
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '9.14'
$c.ClearFormats()
$ws.Range("E29").Value = '  -7.08%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '19.07'
$c.ClearFormats()
$ws.Range("E30").Value = '  -5.12%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '146.22'
$c.ClearFormats()
$ws.Range("E31").Value = '  -5.02%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '5.28'
$c.ClearFormats()
$ws.Range("E32").Value = '  -11.37%  '

# Row 33
$ws.Range("E33").Value = '  -5.43%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.0730'
$c.ClearFormats()
$ws.Range("E34").Value = '  -8.44%  '

# Row 35
$ws.Range("E35").Value = '  -3.54%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.86'
$c.ClearFormats()
$ws.Range("E36").Value = '  -8.56%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.102'
$c.ClearFormats()
$ws.Range("E37").Value = '  -6.04%  '

# Row 38
$ws.Range("E38").Value = '  -11.03%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '3.48'
$c.ClearFormats()
$ws.Range("E39").Value = '  -8.47%  '

# Row 40
$ws.Range("E40").Value = '  -8.65%  '

# Row 41
$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '3.04'
$c.ClearFormats()
$ws.Range("E41").Value = '  -11.89%  '

# Row 42
$ws.Range("E42").Value = '  -0.35%  '

# Row 43
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '13.03'
$c.ClearFormats()
$ws.Range("E43").Value = '  -10.91%  '

# Row 44
$ws.Range("D44").Value = '1.756.70'
$ws.Range("E44").Value = '  +0.69%  '

# Row 45
$ws.Range("E45").Value = '  +2.21%  '

# Row 46
$ws.Range("B46").Value = 'HuobiToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.79'
$c.ClearFormats()
$ws.Range("E46").Value = '  +10.83%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '72.17'
$c.ClearFormats()
$ws.Range("E47").Value = '  -10.28%  '

# Row 48
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '13.74'
$c.ClearFormats()
$ws.Range("E48").Value = '  +5.75%  '

# Row 49
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.171'
$c.ClearFormats()
$ws.Range("E49").Value = '  -11.66%  '

# Row 50
$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '65.22'
$c.ClearFormats()
$ws.Range("E50").Value = '  -8.12%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '91.29'
$c.ClearFormats()
$ws.Range("E51").Value = '  -8.11%  '
